# Weekly update: insert a new price record as row 78 (pushing the
# existing rows 78-102 down to 79-103) on the Hortaliza / Vega
# Monumental Concepción - Alcachofa sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 78..102 down to 79..103, duplicating formatting
# from the row above (matches the style carried on column D).
$ws.Rows.Item(78).Insert()

# Populate the newly-inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value2 = 45135
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 100112013
$ws.Range("G78").Value = "Alcachofa"
$ws.Range("H78").Value = "Española"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 100
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 16000
$ws.Range("M78").Value = 15500
$ws.Range("N78").Value = "$/caja 30 unidades"
$ws.Range("O78").Value = "Provincia de Limarí"
$ws.Range("P78").Value = 517
$ws.Range("Q78").Value = 30
$ws.Range("R78").Value = "Hortaliza"
